$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 206201.218324
$ws.Range("D2").Value = 231.635977

$ws.Range("B3").Value = 28654.054547
$ws.Range("D3").Value = 16.094255
$ws.Range("E3").Value = 0

$ws.Range("B4").Value = 295544.783159
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = -8.678435
$ws.Range("H5").Value = -19.495929
$ws.Range("I5").Value = 2.139059
$ws.Range("J5").Value = 0.143517

$ws.Range("G6").Value = 11.910136
$ws.Range("H6").Value = 0.5481
$ws.Range("I6").Value = 23.272172
$ws.Range("J6").Value = 0.037415

$ws.Range("G7").Value = 20.58857
$ws.Range("H7").Value = 12.045078
$ws.Range("I7").Value = 29.132063
$ws.Range("J7").Value = 0
